$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (timestamp, value, vazao) continuing the existing series,
# one row per hour after the last existing row (128).
$dates = @(45219.361805555556, 45219.40347222222, 45219.445138888892, 45219.486805555556, 45219.52847222222, 45219.570138888892, 45219.611805555556)

$startRow = 129
$endRow = $startRow + $dates.Length - 1

# Copy the formatting of the last existing data row down onto the new rows so
# the new cells pick up the same number formats/styles (date format for col A,
# 2-decimal number format for col B, default for col C) as the rest of the table.
$ws.Range("A128:C128").Copy()
$ws.Range("A" + $startRow + ":C" + $endRow).PasteSpecial(-4122)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 2888.3
}

# Match the saved view state from the diff (scroll position + selection)
$excel.Goto($ws.Range("A123"), $true)
$ws.Range("B126:B135").Select()
